# --------------------------------------------------------------------------
# Dic_Abkuerzungen worksheet update
#
# Four new abbreviation entries are inserted into the alphabetically (desc.)
# ordered list, pushing every following row down by one (cumulatively up to
# four rows by the end of the sheet):
#   - "VN"  / Vereinte Nationen (United Nations)        -> inserted before VGR
#   - "v. a." / vor allem                                -> inserted before usw.
#   - "UN"  / Vereinte Nationen (United Nations)        -> inserted before UK
#   - "Nr." / Nummer / Number                            -> inserted before NOx
#
# Two small text corrections are also carried along with the shift:
#   - BMZ English translation: "Developmen" -> "Development"
#   - EGW German text: "der Gewerblichen Wirtschaft" -> "der Gewerbliche Wirtschaft"
# --------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 5-136 already exist; overwrite their contents to reflect the shifted /
# corrected abbreviation list (keeps existing row formatting).
$ws.Cells.Item(5, 1).Value = "VN"
$ws.Cells.Item(5, 2).Value = "Vereinte Nationen (United Nations)"
$ws.Cells.Item(5, 3).Value = "United Nations"
$ws.Cells.Item(6, 1).Value = "VGR"
$ws.Cells.Item(6, 2).Value = "Volkswirtschaftlichen Gesamtrechnungen"
$ws.Cells.Item(6, 3).Value = ""
$ws.Cells.Item(7, 1).Value = "v. a."
$ws.Cells.Item(7, 2).Value = "vor allem"
$ws.Cells.Item(7, 3).Value = ""
$ws.Cells.Item(8, 1).Value = "usw."
$ws.Cells.Item(8, 2).Value = "und so weiter"
$ws.Cells.Item(8, 3).Value = ""
$ws.Cells.Item(9, 1).Value = "USD"
$ws.Cells.Item(9, 2).Value = "US-Dollar"
$ws.Cells.Item(9, 3).Value = "United States dollar"
$ws.Cells.Item(10, 1).Value = "USA"
$ws.Cells.Item(10, 2).Value = "Vereinigte Staaten von Amerika (United States of America)"
$ws.Cells.Item(10, 3).Value = "United States of America"
$ws.Cells.Item(11, 1).Value = "US"
$ws.Cells.Item(11, 2).Value = "Vereinigte Staaten von Amerika (United States)"
$ws.Cells.Item(11, 3).Value = "United States"
$ws.Cells.Item(12, 1).Value = "UNFCCC"
$ws.Cells.Item(12, 2).Value = ""
$ws.Cells.Item(12, 3).Value = "United Nations Framework Convention on Climate Change"
$ws.Cells.Item(13, 1).Value = "UNCCD"
$ws.Cells.Item(13, 2).Value = "Übereinkommens der Vereinten Nationen zur Bekämpfung der Wüstenbildung in Entwicklungs- und Schwellenländern"
$ws.Cells.Item(13, 3).Value = "UN Convention to Combat Desertification"
$ws.Cells.Item(14, 1).Value = "UN"
$ws.Cells.Item(14, 2).Value = "Vereinte Nationen (United Nations)"
$ws.Cells.Item(14, 3).Value = "United Nations"
$ws.Cells.Item(15, 1).Value = "UK"
$ws.Cells.Item(15, 2).Value = "Vereinigtes Königreich (United Kingdom)"
$ws.Cells.Item(15, 3).Value = "United Kingdom"
$ws.Cells.Item(16, 1).Value = "UBA"
$ws.Cells.Item(16, 2).Value = "Umweltbundesamt"
$ws.Cells.Item(16, 3).Value = "Federal Environment Agency"
$ws.Cells.Item(17, 1).Value = "u.a."
$ws.Cells.Item(17, 2).Value = "unter anderem"
$ws.Cells.Item(17, 3).Value = ""
$ws.Cells.Item(18, 1).Value = "u. a."
$ws.Cells.Item(18, 2).Value = "unter anderem"
$ws.Cells.Item(18, 3).Value = ""
$ws.Cells.Item(19, 1).Value = "TWh"
$ws.Cells.Item(19, 2).Value = "Terawattstunde"
$ws.Cells.Item(19, 3).Value = "Terawatt hour"
$ws.Cells.Item(20, 1).Value = "TREMOD"
$ws.Cells.Item(20, 2).Value = "Transport Emission Estimation Model"
$ws.Cells.Item(20, 3).Value = "Transport Emission Estimation Model"
$ws.Cells.Item(21, 1).Value = "TKU"
$ws.Cells.Item(21, 2).Value = "Telekommunikationsunternehmen"
$ws.Cells.Item(21, 3).Value = "Telecommunications companies"
$ws.Cells.Item(22, 1).Value = "SOEP"
$ws.Cells.Item(22, 2).Value = "Sozio-oekonomischen Panel"
$ws.Cells.Item(22, 3).Value = "Socio-Economic Panel"
$ws.Cells.Item(23, 1).Value = "SO₂"
$ws.Cells.Item(23, 2).Value = "Schwefeldioxid"
$ws.Cells.Item(23, 3).Value = "Sulphur dioxide"
$ws.Cells.Item(24, 1).Value = "SMEs"
$ws.Cells.Item(24, 2).Value = ""
$ws.Cells.Item(24, 3).Value = "Small and medium-sized enterprises"
$ws.Cells.Item(25, 1).Value = "SF₆"
$ws.Cells.Item(25, 2).Value = ""
$ws.Cells.Item(25, 3).Value = "Sulphur hexafluoride"
$ws.Cells.Item(26, 1).Value = "SES"
$ws.Cells.Item(26, 2).Value = "Sozioökonomischer Status"
$ws.Cells.Item(26, 3).Value = "Socioeconomic status"
$ws.Cells.Item(27, 1).Value = "SE"
$ws.Cells.Item(27, 2).Value = "Europäische Gesellschaft"
$ws.Cells.Item(27, 3).Value = ""
$ws.Cells.Item(28, 1).Value = "SDGs"
$ws.Cells.Item(28, 2).Value = "Ziele für Nachhaltige Entwicklung (Sustainable Development Goals)"
$ws.Cells.Item(28, 3).Value = "Sustainable Development Goals"
$ws.Cells.Item(29, 1).Value = "SDG"
$ws.Cells.Item(29, 2).Value = "Ziele für Nachhaltige Entwicklung (Sustainable Development Goals)"
$ws.Cells.Item(29, 3).Value = "Sustainable Development Goals"
$ws.Cells.Item(30, 1).Value = "SALW"
$ws.Cells.Item(30, 2).Value = ""
$ws.Cells.Item(30, 3).Value = "Small arms and light weapons"
$ws.Cells.Item(31, 1).Value = "RKI"
$ws.Cells.Item(31, 2).Value = "Robert Koch-Institut"
$ws.Cells.Item(31, 3).Value = "Robert Koch-Institute"
$ws.Cells.Item(32, 1).Value = "REDD"
$ws.Cells.Item(32, 2).Value = "Verringerung von Emissionen aus Entwaldung und Waldschädigung sowie die Rolle des Waldschutzes, der nachhaltigen Waldbewirtschaftung und des Ausbaus des Kohlenstoffspeichers Wald in Entwicklungsländern"
$ws.Cells.Item(32, 3).Value = "Reducing emissions from deforestation and forest degradation and the role of conservation, sustainable management of forests and enhancement of forest carbon stocks in developing countries"
$ws.Cells.Item(33, 1).Value = "R&D"
$ws.Cells.Item(33, 2).Value = ""
$ws.Cells.Item(33, 3).Value = "Research and development"
$ws.Cells.Item(34, 1).Value = "PM₂.₅"
$ws.Cells.Item(34, 2).Value = "Feinstaub b (Durchmesser kleiner 2,5 Mikrometer)"
$ws.Cells.Item(34, 3).Value = "Particulate matter (diameter smaller than 2.5 micrometers)"
$ws.Cells.Item(35, 1).Value = "PM₂,₅"
$ws.Cells.Item(35, 2).Value = "Feinstaub b (Durchmesser kleiner 2,5 Mikrometer)"
$ws.Cells.Item(35, 3).Value = "Particulate matter (diameter smaller than 2.5 micrometers)"
$ws.Cells.Item(36, 1).Value = "PM₁₀"
$ws.Cells.Item(36, 2).Value = "Feinstaub (Durchmesser kleiner 10 Mikrometer)"
$ws.Cells.Item(36, 3).Value = "Particulate matter (diameter smaller than 10 micrometers)"
$ws.Cells.Item(37, 1).Value = "PM₀.₁"
$ws.Cells.Item(37, 2).Value = "Feinstaub (Durchmesser kleiner 0,1 Mikrometer)"
$ws.Cells.Item(37, 3).Value = "Particulate matter (diameter smaller than 0.1 micrometers)"
$ws.Cells.Item(38, 1).Value = "PM₀,₁"
$ws.Cells.Item(38, 2).Value = "Feinstaub (Durchmesser kleiner 0,1 Mikrometer)"
$ws.Cells.Item(38, 3).Value = "Particulate matter (diameter smaller than 0.1 micrometers)"
$ws.Cells.Item(39, 1).Value = "PKS"
$ws.Cells.Item(39, 2).Value = "Polizeilichen Kriminalstatistik"
$ws.Cells.Item(39, 3).Value = "Police Crime Statistics"
$ws.Cells.Item(40, 1).Value = "PINETI"
$ws.Cells.Item(40, 2).Value = "Pollutant INput and EcosysTem Impact"
$ws.Cells.Item(40, 3).Value = "Pollutant INput and EcosysTem Impact"
$ws.Cells.Item(41, 1).Value = "PFCs"
$ws.Cells.Item(41, 2).Value = ""
$ws.Cells.Item(41, 3).Value = "Perfluorocarbons"
$ws.Cells.Item(42, 1).Value = "P97"
$ws.Cells.Item(42, 2).Value = "97. Perzentil"
$ws.Cells.Item(42, 3).Value = "97th percentile"
$ws.Cells.Item(43, 1).Value = "P90"
$ws.Cells.Item(43, 2).Value = "90. Perzentil"
$ws.Cells.Item(43, 3).Value = "90th percentile"
$ws.Cells.Item(44, 1).Value = "OECD"
$ws.Cells.Item(44, 2).Value = "Organisation für wirtschaftliche Zusammenarbeit und Entwicklung (Organisation for Economic Co-operation and Development)"
$ws.Cells.Item(44, 3).Value = "Organisation for Economic Co-operation and Development"
$ws.Cells.Item(45, 1).Value = "ODA"
$ws.Cells.Item(45, 2).Value = "Öffentliche Entwicklungsausgaben (official development assistance)"
$ws.Cells.Item(45, 3).Value = "Official development assistance"
$ws.Cells.Item(46, 1).Value = "Nr."
$ws.Cells.Item(46, 2).Value = "Nummer"
$ws.Cells.Item(46, 3).Value = "Number"
$ws.Cells.Item(47, 1).Value = "NOₓ"
$ws.Cells.Item(47, 2).Value = "Stickstoffoxid"
$ws.Cells.Item(47, 3).Value = "Nitrogen oxides"
$ws.Cells.Item(48, 1).Value = "NMVOCs"
$ws.Cells.Item(48, 2).Value = "Flüchtige organische Verbindungen (non-methane volatile organic compounds)"
$ws.Cells.Item(48, 3).Value = "Non-methane volatile organic compounds"
$ws.Cells.Item(49, 1).Value = "NMVOC"
$ws.Cells.Item(49, 2).Value = "Flüchtige organische Verbindungen (non-methane volatile organic compounds)"
$ws.Cells.Item(49, 3).Value = "non-methane volatile organic compounds"
$ws.Cells.Item(50, 1).Value = "NH₃"
$ws.Cells.Item(50, 2).Value = "Ammoniak"
$ws.Cells.Item(50, 3).Value = "Ammonia"
$ws.Cells.Item(51, 1).Value = "NF₃"
$ws.Cells.Item(51, 2).Value = "Stickstofftrifluorid"
$ws.Cells.Item(51, 3).Value = "Nitrogen trifluoride"
$ws.Cells.Item(52, 1).Value = "NEC"
$ws.Cells.Item(52, 2).Value = "Richtlinie über nationale Emissionshöchstmengen für bestimmte Luftschadstoffe (National Emission Ceilings Directive)"
$ws.Cells.Item(52, 3).Value = "National Emission Ceilings Directive"
$ws.Cells.Item(53, 1).Value = "N₂O"
$ws.Cells.Item(53, 2).Value = "Lachgas"
$ws.Cells.Item(53, 3).Value = "Nitrous oxide"
$ws.Cells.Item(54, 1).Value = "N"
$ws.Cells.Item(54, 2).Value = ""
$ws.Cells.Item(54, 3).Value = "Nitrogen"
$ws.Cells.Item(55, 1).Value = "MSY"
$ws.Cells.Item(55, 2).Value = "Maximum Sustainable Yield"
$ws.Cells.Item(55, 3).Value = "Maximum Sustainable Yield"
$ws.Cells.Item(56, 1).Value = "Mrd."
$ws.Cells.Item(56, 2).Value = "Milliarde"
$ws.Cells.Item(56, 3).Value = ""
$ws.Cells.Item(57, 1).Value = "mg/l"
$ws.Cells.Item(57, 2).Value = "Milligramm pro Liter"
$ws.Cells.Item(57, 3).Value = "Miligrams per litre"
$ws.Cells.Item(58, 1).Value = "mg"
$ws.Cells.Item(58, 2).Value = "Milligramm"
$ws.Cells.Item(58, 3).Value = "Miligrams"
$ws.Cells.Item(59, 1).Value = "Mbps"
$ws.Cells.Item(59, 2).Value = ""
$ws.Cells.Item(59, 3).Value = "Megabit per second"
$ws.Cells.Item(60, 1).Value = "Mbit/s"
$ws.Cells.Item(60, 2).Value = "Megabit pro Sekunde"
$ws.Cells.Item(60, 3).Value = "Megabit per second"
$ws.Cells.Item(61, 1).Value = "m³"
$ws.Cells.Item(61, 2).Value = "Kubikmeter"
$ws.Cells.Item(61, 3).Value = "Cubic metre"
$ws.Cells.Item(62, 1).Value = "m²"
$ws.Cells.Item(62, 2).Value = "Quadratmeter"
$ws.Cells.Item(62, 3).Value = "Square meter"
$ws.Cells.Item(63, 1).Value = "LULUCF"
$ws.Cells.Item(63, 2).Value = ""
$ws.Cells.Item(63, 3).Value = "Land use, land-use change and forestry"
$ws.Cells.Item(64, 1).Value = "LDCs"
$ws.Cells.Item(64, 2).Value = "am wenigsten entwickelte Länder (Least developed countries)"
$ws.Cells.Item(64, 3).Value = "Least developed countries"
$ws.Cells.Item(65, 1).Value = "LDC"
$ws.Cells.Item(65, 2).Value = "am wenigsten entwickelte Länder (Least developed countries)"
$ws.Cells.Item(65, 3).Value = "Least developed countries"
$ws.Cells.Item(66, 1).Value = "LAWA"
$ws.Cells.Item(66, 2).Value = "Bund/Länder-Arbeitsgemeinschaft Wasser"
$ws.Cells.Item(66, 3).Value = "German Working Group on Water Issues of the Länder and the Federal Government"
$ws.Cells.Item(67, 1).Value = "l"
$ws.Cells.Item(67, 2).Value = "Liter"
$ws.Cells.Item(67, 3).Value = "Litre"
$ws.Cells.Item(68, 1).Value = "km²"
$ws.Cells.Item(68, 2).Value = "Quadratkilometer"
$ws.Cells.Item(68, 3).Value = "Square kilometer"
$ws.Cells.Item(69, 1).Value = "KiGGS"
$ws.Cells.Item(69, 2).Value = "Studie zur Gesundheit von Kindern und Jugendlichen in Deutschland"
$ws.Cells.Item(69, 3).Value = "Study on the health of children and adolescents in Germany"
$ws.Cells.Item(70, 1).Value = "kg/m²"
$ws.Cells.Item(70, 2).Value = "Kilogramm pro Quadratmeter"
$ws.Cells.Item(70, 3).Value = "Kilogram per square meter"
$ws.Cells.Item(71, 1).Value = "kg/ha"
$ws.Cells.Item(71, 2).Value = "Kilogramm pro Hektar"
$ws.Cells.Item(71, 3).Value = "Kilogram per hectare"
$ws.Cells.Item(72, 1).Value = "kg"
$ws.Cells.Item(72, 2).Value = "Kilogramm"
$ws.Cells.Item(72, 3).Value = "Kilogram"
$ws.Cells.Item(73, 1).Value = "Kfz"
$ws.Cells.Item(73, 2).Value = "Kraftfahrzeug"
$ws.Cells.Item(73, 3).Value = ""
$ws.Cells.Item(74, 1).Value = "KfW"
$ws.Cells.Item(74, 2).Value = "Kreditanstalt für Wiederaufbau"
$ws.Cells.Item(74, 3).Value = "Kreditanstalt für Wiederaufbau"
$ws.Cells.Item(75, 1).Value = "ISCO"
$ws.Cells.Item(75, 2).Value = "Internationale Standardklassifikation der Berufe (International Standard Classification of Occupations)"
$ws.Cells.Item(75, 3).Value = "International Standard Classification of Occupations"
$ws.Cells.Item(76, 1).Value = "ISCED"
$ws.Cells.Item(76, 2).Value = "International Standard Classification of Education"
$ws.Cells.Item(76, 3).Value = "International Standard Classification of Education"
$ws.Cells.Item(77, 1).Value = "i.e."
$ws.Cells.Item(77, 2).Value = ""
$ws.Cells.Item(77, 3).Value = "that is to say (id est)"
$ws.Cells.Item(78, 1).Value = "H-FKW/HFC"
$ws.Cells.Item(78, 2).Value = "Teilhalogenierte Fluorkohlenwasserstoffe"
$ws.Cells.Item(78, 3).Value = ""
$ws.Cells.Item(79, 1).Value = "HFCs"
$ws.Cells.Item(79, 2).Value = ""
$ws.Cells.Item(79, 3).Value = "Hydrofluorocarbons"
$ws.Cells.Item(80, 1).Value = "ha"
$ws.Cells.Item(80, 2).Value = "Hektar"
$ws.Cells.Item(80, 3).Value = "Hectare"
$ws.Cells.Item(81, 1).Value = "GPG"
$ws.Cells.Item(81, 2).Value = "Geschlechtsspezifischen Verdienstabstand (gender pay gap)"
$ws.Cells.Item(81, 3).Value = "gender pay gap"
$ws.Cells.Item(82, 1).Value = "GNI"
$ws.Cells.Item(82, 2).Value = ""
$ws.Cells.Item(82, 3).Value = "Gross national income"
$ws.Cells.Item(83, 1).Value = "GmbH"
$ws.Cells.Item(83, 2).Value = "Gesellschaft mit beschränkter Haftung"
$ws.Cells.Item(83, 3).Value = "Company with limited liability"
$ws.Cells.Item(84, 1).Value = "GIZ"
$ws.Cells.Item(84, 2).Value = "Deutsche Gesellschaft für Internationale Zusammenarbeit"
$ws.Cells.Item(84, 3).Value = "Deutsche Gesellschaft für Internationale Zusammenarbeit"
$ws.Cells.Item(85, 1).Value = "gGmbH"
$ws.Cells.Item(85, 2).Value = "gemeinnützige Gesellschaft mit beschränkter Haftung"
$ws.Cells.Item(85, 3).Value = "Non-profit limited liability company"
$ws.Cells.Item(86, 1).Value = "GG"
$ws.Cells.Item(86, 2).Value = "Grundgesetz"
$ws.Cells.Item(86, 3).Value = "Basic Law"
$ws.Cells.Item(87, 1).Value = "gender pay gap"
$ws.Cells.Item(87, 2).Value = "Geschlechtsspezifischen Verdienstabstand"
$ws.Cells.Item(87, 3).Value = ""
$ws.Cells.Item(88, 1).Value = "GDP"
$ws.Cells.Item(88, 2).Value = ""
$ws.Cells.Item(88, 3).Value = "Gross domestic product"
$ws.Cells.Item(89, 1).Value = "FuE"
$ws.Cells.Item(89, 2).Value = "Forschung und Entwicklung"
$ws.Cells.Item(89, 3).Value = ""
$ws.Cells.Item(90, 1).Value = "FTTB/H"
$ws.Cells.Item(90, 2).Value = "Reine Glasfasernetze"
$ws.Cells.Item(90, 3).Value = "Fully fibre-optic networks"
$ws.Cells.Item(91, 1).Value = "FKW/PFC"
$ws.Cells.Item(91, 2).Value = "Perfluorierte Kohlenwasserstoffe"
$ws.Cells.Item(91, 3).Value = ""
$ws.Cells.Item(92, 1).Value = "FidAR"
$ws.Cells.Item(92, 2).Value = "Frauen in die Aufsichtsräte"
$ws.Cells.Item(92, 3).Value = ""
$ws.Cells.Item(93, 1).Value = "FCPF"
$ws.Cells.Item(93, 2).Value = "Forest Carbon Partnership Facility"
$ws.Cells.Item(93, 3).Value = "Forest Carbon Partnership Facility"
$ws.Cells.Item(94, 1).Value = "FAO"
$ws.Cells.Item(94, 2).Value = "Ernährungs- und Landwirtschaftsorganisation der Vereinten Nationen (Food and Agriculture Organization)"
$ws.Cells.Item(94, 3).Value = "Food and Agriculture Organization"
$ws.Cells.Item(95, 1).Value = "EU-SILC"
$ws.Cells.Item(95, 2).Value = "Statistik über Einkommen und Lebensbedingungen (Statistics on Income and Living Conditions)"
$ws.Cells.Item(95, 3).Value = "Statistics on Income and Living Conditions"
$ws.Cells.Item(96, 1).Value = "EUR"
$ws.Cells.Item(96, 2).Value = "Euro"
$ws.Cells.Item(96, 3).Value = "Euro"
$ws.Cells.Item(97, 1).Value = "EU-EVK"
$ws.Cells.Item(97, 2).Value = "EU-Energieverbrauchskennzeichnung"
$ws.Cells.Item(97, 3).Value = ""
$ws.Cells.Item(98, 1).Value = "EUA"
$ws.Cells.Item(98, 2).Value = "Europäische Umweltagentur"
$ws.Cells.Item(98, 3).Value = ""
$ws.Cells.Item(99, 1).Value = "EU-28"
$ws.Cells.Item(99, 2).Value = "Europäische Union mit 28 Mitgliedsstaaten"
$ws.Cells.Item(99, 3).Value = "European Union consisting of 28 member states"
$ws.Cells.Item(100, 1).Value = "EU-27"
$ws.Cells.Item(100, 2).Value = "Europäische Union mit 27 Mitgliedsstaaten"
$ws.Cells.Item(100, 3).Value = "European Union consisting of 27 member states"
$ws.Cells.Item(101, 1).Value = "EU"
$ws.Cells.Item(101, 2).Value = "Europäische Union"
$ws.Cells.Item(101, 3).Value = "European Union"
$ws.Cells.Item(102, 1).Value = "etc."
$ws.Cells.Item(102, 2).Value = "und so weiter (et cetera)"
$ws.Cells.Item(102, 3).Value = "and so on (et cetera)"
$ws.Cells.Item(103, 1).Value = "ESVG"
$ws.Cells.Item(103, 2).Value = "Europäische System Volkswirtschaftlicher Gesamtrechnungen"
$ws.Cells.Item(103, 3).Value = ""
$ws.Cells.Item(104, 1).Value = "ESA"
$ws.Cells.Item(104, 2).Value = ""
$ws.Cells.Item(104, 3).Value = "European System of National and Regional Accounts"
$ws.Cells.Item(105, 1).Value = "EMAS"
$ws.Cells.Item(105, 2).Value = "Eco-Management and Audit Scheme"
$ws.Cells.Item(105, 3).Value = "Eco-Management and Audit Scheme"
$ws.Cells.Item(106, 1).Value = "EGW"
$ws.Cells.Item(106, 2).Value = "Ernährungs- und der Gewerbliche Wirtschaft"
$ws.Cells.Item(106, 3).Value = "Food and industrial economy"
$ws.Cells.Item(107, 1).Value = "EEG"
$ws.Cells.Item(107, 2).Value = "Erneuerbare-Energien-Gesetz"
$ws.Cells.Item(107, 3).Value = "Renewable Energy Sources Act"
$ws.Cells.Item(108, 1).Value = "EEA"
$ws.Cells.Item(108, 2).Value = ""
$ws.Cells.Item(108, 3).Value = "European Environment Agency"
$ws.Cells.Item(109, 1).Value = "e.g."
$ws.Cells.Item(109, 2).Value = ""
$ws.Cells.Item(109, 3).Value = "for example (exempli gratia)"
$ws.Cells.Item(110, 1).Value = "DIN"
$ws.Cells.Item(110, 2).Value = "Deutsches Institut für Normung e.V."
$ws.Cells.Item(110, 3).Value = "German Institute for Standardisation Registered Association"
$ws.Cells.Item(111, 1).Value = "DEG"
$ws.Cells.Item(111, 2).Value = "Deutsche Investitions- und Entwicklungsgesellschaft"
$ws.Cells.Item(111, 3).Value = "Deutsche Investitions- und Entwicklungsgesellschaft"
$ws.Cells.Item(112, 1).Value = "DDB"
$ws.Cells.Item(112, 2).Value = "Deutsche Digitale Bibliothek"
$ws.Cells.Item(112, 3).Value = "German Digital Library (Deutsche Digitale Bibliothek)"
$ws.Cells.Item(113, 1).Value = "DDA"
$ws.Cells.Item(113, 2).Value = "Dachverband Deutscher Avifaunisten"
$ws.Cells.Item(113, 3).Value = ""
$ws.Cells.Item(114, 1).Value = "DAC"
$ws.Cells.Item(114, 2).Value = "Richtlinien des Entwicklungsausschusses (Development Assistance Committee)"
$ws.Cells.Item(114, 3).Value = "Development Assistance Committee"
$ws.Cells.Item(115, 1).Value = "CPI"
$ws.Cells.Item(115, 2).Value = "Corruption Perception Index"
$ws.Cells.Item(115, 3).Value = "Corruption Perception Index"
$ws.Cells.Item(116, 1).Value = "COVID-19"
$ws.Cells.Item(116, 2).Value = "Coronavirus SARS-CoV-2"
$ws.Cells.Item(116, 3).Value = "Coronavirus SARS-CoV-2"
$ws.Cells.Item(117, 1).Value = "CO₂"
$ws.Cells.Item(117, 2).Value = "Kohlenstoffdioxid"
$ws.Cells.Item(117, 3).Value = "Carbon dioxide"
$ws.Cells.Item(118, 1).Value = "CLRTAP"
$ws.Cells.Item(118, 2).Value = "Genfer Luftreinhaltekonvention (Convention on Long-Range Transboundary Air Pollution)"
$ws.Cells.Item(118, 3).Value = "Convention on Long-Range Transboundary Air Pollution"
$ws.Cells.Item(119, 1).Value = "CH₄"
$ws.Cells.Item(119, 2).Value = "Methan"
$ws.Cells.Item(119, 3).Value = "Methane"
$ws.Cells.Item(120, 1).Value = "CATV"
$ws.Cells.Item(120, 2).Value = "Kabelfernsehen"
$ws.Cells.Item(120, 3).Value = "Cable television"
$ws.Cells.Item(121, 1).Value = "bzw."
$ws.Cells.Item(121, 2).Value = "beziehungsweise"
$ws.Cells.Item(121, 3).Value = ""
$ws.Cells.Item(122, 1).Value = "BNE"
$ws.Cells.Item(122, 2).Value = "Bruttonationaleinkommen"
$ws.Cells.Item(122, 3).Value = ""
$ws.Cells.Item(123, 1).Value = "bn"
$ws.Cells.Item(123, 2).Value = ""
$ws.Cells.Item(123, 3).Value = "Billion"
$ws.Cells.Item(124, 1).Value = "BMZ"
$ws.Cells.Item(124, 2).Value = "Bundesministerium für wirtschaftliche Zusammenarbeit und Entwicklung"
$ws.Cells.Item(124, 3).Value = "Federal Ministry for Economic Cooperation and Development"
$ws.Cells.Item(125, 1).Value = "BMVI"
$ws.Cells.Item(125, 2).Value = "Bundesministerium für Verkehr und digitale Infrastruktur"
$ws.Cells.Item(125, 3).Value = "Federal Ministry of Transport and Digital Infrastructure"
$ws.Cells.Item(126, 1).Value = "BMI"
$ws.Cells.Item(126, 2).Value = "Body Mass Index"
$ws.Cells.Item(126, 3).Value = "Body Mass Index"
$ws.Cells.Item(127, 1).Value = "BMEL"
$ws.Cells.Item(127, 2).Value = "Bundesministeriums für Ernährung und Landwirtschaft"
$ws.Cells.Item(127, 3).Value = "Federal Ministry of Food and Agriculture"
$ws.Cells.Item(128, 1).Value = "BLE"
$ws.Cells.Item(128, 2).Value = "Bundesanstalt für Landwirtschaft und Ernährung"
$ws.Cells.Item(128, 3).Value = "Federal Office for Agriculture and Food"
$ws.Cells.Item(129, 1).Value = "BKG"
$ws.Cells.Item(129, 2).Value = "Bundesamt für Kartographie und Geodäsie"
$ws.Cells.Item(129, 3).Value = "Federal Agency for Cartography and Geodesy"
$ws.Cells.Item(130, 1).Value = "BIP"
$ws.Cells.Item(130, 2).Value = "Bruttoinlandsprodukt"
$ws.Cells.Item(130, 3).Value = ""
$ws.Cells.Item(131, 1).Value = "BfN"
$ws.Cells.Item(131, 2).Value = "Bundesamt für Naturschutz"
$ws.Cells.Item(131, 3).Value = ""
$ws.Cells.Item(132, 1).Value = "BEEG"
$ws.Cells.Item(132, 2).Value = "Bundeselterngeld- und Elternzeitgesetz"
$ws.Cells.Item(132, 3).Value = "Federal Parental Allowance and Parental Leave Act"
$ws.Cells.Item(133, 1).Value = "Art."
$ws.Cells.Item(133, 2).Value = "Artikel"
$ws.Cells.Item(133, 3).Value = "Article"
$ws.Cells.Item(134, 1).Value = "ALKIS"
$ws.Cells.Item(134, 2).Value = "Amtlichen Liegenschaftskataster-Informationssystem"
$ws.Cells.Item(134, 3).Value = "Official land register information system"
$ws.Cells.Item(135, 1).Value = "AGEE-Stat"
$ws.Cells.Item(135, 2).Value = "Arbeitsgruppe Erneuerbare Energien-Statistik"
$ws.Cells.Item(135, 3).Value = "Working Group on Renewable Energy Statistics"
$ws.Cells.Item(136, 1).Value = "AGEB"
$ws.Cells.Item(136, 2).Value = "Arbeitsgemeinschaft Energiebilanzen"
$ws.Cells.Item(136, 3).Value = "Energy Balance Association"

# Rows 137-140 are brand new; write their values first, then copy the cell
# formatting from the last pre-existing data row (136) so the new rows match
# the look of the rest of the table.
$ws.Cells.Item(137, 1).Value = "Abs."
$ws.Cells.Item(137, 2).Value = "Absatz"
$ws.Cells.Item(137, 3).Value = ""
$ws.Cells.Item(138, 1).Value = "a.m."
$ws.Cells.Item(138, 2).Value = ""
$ws.Cells.Item(138, 3).Value = "before noon (ante meridiem)"
$ws.Cells.Item(139, 1).Value = "µg/m³"
$ws.Cells.Item(139, 2).Value = "Mikrogramm pro Kubikmeter"
$ws.Cells.Item(139, 3).Value = "Micrograms per cubic metre"
$ws.Cells.Item(140, 1).Value = "µg"
$ws.Cells.Item(140, 2).Value = "Mikrogramm"
$ws.Cells.Item(140, 3).Value = "Micrograms"

$formatSource = $ws.Range("A136:C136")
$formatTarget = $ws.Range("A137:C140")
$formatSource.Copy() | Out-Null
$formatTarget.PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Keep the explicit selection in sync with the original workbook.
$ws.Range("C5").Select() | Out-Null
